# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" list (column E, rows 16-27) is reordered from
# descending (2104 .. 2005) to ascending (2005 .. 2104), and the
# "Valor Mora" (column F) values follow the same period so that the
# distinct value (31200) stays attached to period 2104, which is now
# the last row instead of the first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$periods = @{
    16 = "2005"
    17 = "2006"
    18 = "2007"
    19 = "2008"
    20 = "2009"
    21 = "2010"
    22 = "2011"
    23 = "2012"
    24 = "2101"
    25 = "2102"
    26 = "2103"
    27 = "2104"
}

$valores = @{
    16 = 36000
    17 = 36000
    18 = 36000
    19 = 36000
    20 = 36000
    21 = 36000
    22 = 36000
    23 = 36000
    24 = 36000
    25 = 36000
    26 = 36000
    27 = 31200
}

foreach ($row in 16..27) {
    $ws.Range("E$row").Value = $periods[$row]
    $ws.Range("F$row").Value = $valores[$row]
}
